$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RegistrationData")

# Fill in the new "Blank Password" scenario row (row 8)
$ws.Range("A8").Value = "Lokesh"
$ws.Range("B8").Value = "Sharma"
$ws.Range("C8").Value = "lokesh407@xtivia.com"
$ws.Range("D8").Value = 8447520166
$ws.Range("G8").Value = "Invalid Data"
$ws.Range("H8").Value = "Password cant be blank"

# Add hyperlinks for the email and contact cells, matching the pattern used by row 7
$ws.Hyperlinks.Add($ws.Range("C8"), "mailto:lokesh407@xtivia.com")
$ws.Hyperlinks.Add($ws.Range("D8"), "mailto:lokesh403@xtivia.com", "", "", "lokesh403@xtivia.com")
$ws.Range("D8").Value = 8447520166

# Re-apply the existing "hyperlink" formatting style from row 7 so that the
# styling matches the rest of the data (Excel applies a brand-new style object
# when Hyperlinks.Add is used, so copy the formats to reuse the existing one).
$ws.Range("C7:D7").Copy()
$ws.Range("C8:D8").PasteSpecial(-4122)  # xlPasteFormats

# Update the active selection to A8
$ws.Activate()
$ws.Range("A8").Select()
